$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 14.4
$ws.Range("B3").Value = 19.3
$ws.Range("C3").Value = 18.5
$ws.Range("C10").Value = 14.9
$ws.Range("C13").Value = 15.2
$ws.Range("C18").Value = 14.8
